$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "PROJECT 2" table (rows 10-16), added below the existing table ---
# Title row
$ws.Range("A10").Value = "PROJECT 2"

# Write "Register" / "Login" first so they land in the shared-strings table
# ahead of "Server Receive" / "Server Send" / "Authentication" (matches the
# order new strings were interned in the saved workbook).
$ws.Range("C12").Value = "Register"
$ws.Range("C13").Value = "Login"

# Header row for the new table
$ws.Range("A11").Value = "Command ID"
$ws.Range("C11").Value = "Client"
$ws.Range("D11").Value = "Server Receive"
$ws.Range("E11").Value = "Server Send"
$ws.Range("F11").Value = "Authentication"

# Row 12 - Register
$ws.Range("A12").Value = 0
$ws.Range("D12").Value = "Register"
$ws.Range("E12").Value = "Register"

# Row 13 - Login
$ws.Range("A13").Value = 1
$ws.Range("D13").Value = "Login"
$ws.Range("E13").Value = "Login"

# Row 14 - Join Room
$ws.Range("A14").Value = 2
$ws.Range("C14").Value = "Join Room"
$ws.Range("D14").Value = "Join Room"

# Row 15 - Leave Room
$ws.Range("A15").Value = 3
$ws.Range("C15").Value = "Leave Room"
$ws.Range("D15").Value = "Leave Room"

# Row 16 - Message to/others in Room
$ws.Range("A16").Value = 4
$ws.Range("C16").Value = "Message to Room"
$ws.Range("D16").Value = "Message others in Room"

# New columns E/F need custom widths to fit their (longer) header text,
# matching the auto-sized columns Excel produced for the new table.
$ws.Columns.Item(5).ColumnWidth = 13.67
$ws.Columns.Item(6).ColumnWidth = 12.33

# Final selection lands on F12, matching the edited workbook's cursor.
$ws.Range("F12").Select()
